$wb = $excel.ActiveWorkbook

# Sheet "展览" - update 想去人数 (F column) for rows 2, 5, 6
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 3357
$ws1.Range("F5").Value = 1475
$ws1.Range("F6").Value = 37

# Sheet "全部类型" - same events, same updates
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 3357
$ws4.Range("F5").Value = 1475
$ws4.Range("F6").Value = 37
